$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Enterprises density (per 1000 people)
$ws.Cells.Item(11, 2).Value = "'25.71"
$ws.Cells.Item(11, 3).Value = "'0.88"
$ws.Cells.Item(11, 4).Value = "'26.59"

# Row 12: Employment (% of total)
$ws.Cells.Item(12, 2).Value = "'57.38"
$ws.Cells.Item(12, 3).Value = "'17.26"
$ws.Cells.Item(12, 4).Value = "'74.64"

# Row 14: Enterprises (% of total)
$ws.Cells.Item(14, 2).Value = "'96.52"
$ws.Cells.Item(14, 4).Value = "'99.82"
